# Updated ITA_grids model - 2025-08-20 18:31
$wb = $excel.ActiveWorkbook

# --- Sheet "ev_charging_uc": update the two comma-separated timeslice lists.
# C13/C14 are source values; G7/G8 hold formulas (=C14 / =C13) that pick
# these values up automatically on recalculation.
$wsUc = $wb.Worksheets.Item("ev_charging_uc")
$wsUc.Range("C13").Value2 = "WaD,WaP,FaD,FaP,SaP,SaD,RaD,RaP"
$wsUc.Range("C14").Value2 = "SaN,WaN,RaN,RaP,WaP,FaP,SaP,FaN"

# --- Sheet "re_profiles": the M/N columns (season letter + share) for rows
# 4-7 were cyclically re-ordered: old row4->row6, row5->row4, row6->row7,
# row7->row5.
$wsRe = $wb.Worksheets.Item("re_profiles")

$m4 = $wsRe.Range("M4").Value2
$n4 = $wsRe.Range("N4").Value2
$m5 = $wsRe.Range("M5").Value2
$n5 = $wsRe.Range("N5").Value2
$m6 = $wsRe.Range("M6").Value2
$n6 = $wsRe.Range("N6").Value2
$m7 = $wsRe.Range("M7").Value2
$n7 = $wsRe.Range("N7").Value2

$wsRe.Range("M4").Value2 = $m5
$wsRe.Range("N4").Value2 = $n5
$wsRe.Range("M5").Value2 = $m7
$wsRe.Range("N5").Value2 = $n7
$wsRe.Range("M6").Value2 = $m4
$wsRe.Range("N6").Value2 = $n4
$wsRe.Range("M7").Value2 = $m6
$wsRe.Range("N7").Value2 = $n6
